# Auto-generated Excel COM-interop script to apply crypto price/volume updates
# (from 'Updated cryptos list' GitHub Actions commit)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text would NOT be mis-parsed as a number by Excel:
# set them directly so the underlying XML stays a plain (unstyled) string cell.
$ws.Range('D2').Value = '41.315.96'
$ws.Range('E2').Value = '  -3.44%  '
$ws.Range('D3').Value = '2.464.32'
$ws.Range('E3').Value = '  -2.52%  '
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('E5').Value = '  +0.29%  '
$ws.Range('E6').Value = '  -6.33%  '
$ws.Range('E7').Value = '  -2.49%  '
$ws.Range('E8').Value = '  -0.03%  '
$ws.Range('E9').Value = '  -4.57%  '
$ws.Range('E10').Value = '  -6.26%  '
$ws.Range('E11').Value = '  -3.07%  '
$ws.Range('E12').Value = '  -1.05%  '
$ws.Range('E13').Value = '  -4.53%  '
$ws.Range('D14').Value = '2.843.19'
$ws.Range('E14').Value = '  -2.64%  '
$ws.Range('D15').Value = '2.462.86'
$ws.Range('E15').Value = '  -3.14%  '
$ws.Range('E16').Value = '  -3.37%  '
$ws.Range('E17').Value = '  -3.47%  '
$ws.Range('D18').Value = '41.267.50'
$ws.Range('E18').Value = '  -3.51%  '
$ws.Range('E19').Value = '  -5.33%  '
$ws.Range('E20').Value = '  -3.18%  '
$ws.Range('E21').Value = '  -8.98%  '
$ws.Range('E22').Value = '  -1.90%  '
$ws.Range('E23').Value = '  -2.53%  '
$ws.Range('E24').Value = '  -4.15%  '
$ws.Range('E25').Value = '  +0.08%  '
$ws.Range('E26').Value = '  -6.14%  '
$ws.Range('E27').Value = '  -5.27%  '
$ws.Range('E28').Value = '  -5.84%  '
$ws.Range('E29').Value = '  -5.41%  '
$ws.Range('E30').Value = '  -4.72%  '
$ws.Range('E31').Value = '  -6.16%  '
$ws.Range('E32').Value = '  -5.16%  '
$ws.Range('E33').Value = '  -5.41%  '
$ws.Range('E34').Value = '  -2.96%  '
$ws.Range('E35').Value = '  -5.53%  '
$ws.Range('E36').Value = '  -0.93%  '
$ws.Range('E37').Value = '  -3.57%  '
$ws.Range('E38').Value = '  -7.13%  '
$ws.Range('E39').Value = '  -2.83%  '
$ws.Range('B40').Value = 'RenderToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('E40').Value = '  +2.82%  '
$ws.Range('B41').Value = 'Kaspa'
$ws.Range('C41').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('E41').Value = '  -7.82%  '
$ws.Range('E42').Value = '  +0.10%  '
$ws.Range('E43').Value = '  -10.08%  '
$ws.Range('D44').Value = '1.987.75'
$ws.Range('E44').Value = '  -0.03%  '
$ws.Range('E45').Value = '  -4.77%  '
$ws.Range('E46').Value = '  -7.84%  '
$ws.Range('E47').Value = '  -3.12%  '
$ws.Range('E48').Value = '  -4.31%  '
$ws.Range('E49').Value = '  -3.82%  '
$ws.Range('B50').Value = 'BitcoinSV'
$ws.Range('C50').Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range('E50').Value = '  -6.01%  '
$ws.Range('B51').Value = 'Algorand'
$ws.Range('C51').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('E51').Value = '  -6.80%  '

# Cells whose new text looks like a plain decimal number (e.g. '94.43').
# Excel would silently convert these to a Number type and lose the exact
# text/formatting (trailing zeros, etc.), so force a Text number format
# first, assign the literal string, then restore the default 'Normal' style
# so no residual style index is left on the cell.
$textCells = @('D5', 'D6', 'D9', 'D10', 'D16', 'D19', 'D21', 'D22', 'D23', 'D26', 'D27', 'D28', 'D29', 'D30', 'D31', 'D35', 'D36', 'D37', 'D38', 'D40', 'D41', 'D43', 'D45', 'D46', 'D47', 'D48', 'D49', 'D50', 'D51')
foreach ($cell in $textCells) { $ws.Range($cell).NumberFormat = '@' }

$ws.Range('D5').Value = '311.63'
$ws.Range('D6').Value = '94.43'
$ws.Range('D9').Value = '0.499'
$ws.Range('D10').Value = '33.52'
$ws.Range('D16').Value = '14.83'
$ws.Range('D19').Value = '6.32'
$ws.Range('D21').Value = '11.28'
$ws.Range('D22').Value = '68.45'
$ws.Range('D23').Value = '237.25'
$ws.Range('D26').Value = '1.91'
$ws.Range('D27').Value = '24.12'
$ws.Range('D28').Value = '2.21'
$ws.Range('D29').Value = '9.62'
$ws.Range('D30').Value = '36.83'
$ws.Range('D31').Value = '152.06'
$ws.Range('D35').Value = '0.0745'
$ws.Range('D36').Value = '3.05'
$ws.Range('D37').Value = '1.89'
$ws.Range('D38').Value = '17.04'
$ws.Range('D40').Value = '4.29'
$ws.Range('D41').Value = '0.102'
$ws.Range('D43').Value = '19.95'
$ws.Range('D45').Value = '0.0285'
$ws.Range('D46').Value = '3.05'
$ws.Range('D47').Value = '8.84'
$ws.Range('D48').Value = '69.26'
$ws.Range('D49').Value = '96.95'
$ws.Range('D50').Value = '74.74'
$ws.Range('D51').Value = '0.178'

foreach ($cell in $textCells) { $ws.Range($cell).Style = 'Normal' }
